# "Generate Report for Handback" — marks the two localization files (for
# zh-cn and de-de) as handed back and in sync with en-US, stamping the
# Latest Target File / Latest Handback File / Latest Handback DateTime
# columns on each language sheet, and updating the rollup Status on the
# Overview sheet.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: Status is mirrored into columns E (zh-cn) and F (de-de)
# for both tracked files.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de): for each of the two tracked source
# files, record the target file (a link back to the source .md), the
# handback file (same name as the most recent handoff file, since the
# translation is in sync) and the handback timestamp, then refresh Status.
# ---------------------------------------------------------------------
$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8fc6a8797557e79e10d9b08ba8be7e20595628e/e2e"

$files = @(
    @{ Name = "2752bca6-a840-4b1d-980c-31ee1e8b9ea4.md"; Row = 2 },
    @{ Name = "ad1ecfcc-9247-45bd-88b9-6231b4fabed8.md"; Row = 3 }
)

$langs = @(
    @{ Sheet = "zh-cn"; HandbackTime = "2016-11-29 05:13:03" },
    @{ Sheet = "de-de"; HandbackTime = "2016-11-29 05:13:21" }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    $ws.Hyperlinks.Delete()

    foreach ($file in $files) {
        $row = $file.Row
        $url = "$baseUrl/$($file.Name)"

        # Column A: Source File Name (existing hyperlink, re-created so the
        # hyperlink collection keeps A-before-I ordering per row).
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 1), $url, "", "", $file.Name)

        # Column C: Status
        $ws.Cells.Item($row, 3).Value = $statusText

        # Column I: Latest Target File — points at the same source doc.
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 9), $url, "", "", $file.Name)

        # Column J: Latest Handback File — same payload as the latest
        # handoff file (column G) since the translation is in sync.
        $handoffFile = $ws.Cells.Item($row, 7).Value2
        $ws.Cells.Item($row, 10).Value = $handoffFile

        # Column K: Latest Handback DateTime
        $ws.Cells.Item($row, 11).Value = $lang.HandbackTime
    }

    $ws.Columns.Item(3).ColumnWidth = 29.14
    $ws.Columns.Item(9).ColumnWidth = 39.17
    $ws.Columns.Item(10).ColumnWidth = 39.17
}
